$d = $word.ActiveDocument

# Delete the paragraph "Ver no Jupiter Salvar em pdf Salvar em docx" entirely
# (including its own paragraph mark) and the following paragraph
# "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages.
# Original theme under Creative Commons Attribution", also including its
# paragraph mark. Leave the blank paragraph right before them intact.
#
# Walk the paragraphs back-to-front so that deleting one paragraph's range
# does not shift/ invalidate the indices of paragraphs still to be removed.

$paras = @($d.Paragraphs)
for ($i = $paras.Count - 1; $i -ge 0; $i--) {
    $t = $paras[$i].Range.Text
    if ($t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*" -or
        $t -like "*Contact: luizeleno@usp.br*") {
        $paras[$i].Range.Delete()
    }
}
